# Update dilution/weight figures for the 0.02M PBS preparation protocol.
$wb = $excel.ActiveWorkbook

# --- Midazolam sheet ---
$wsMidazolam = $wb.Worksheets.Item("Midazolam")
$wsMidazolam.Range("D9").Value = 0.04

# --- Ketamine sheet ---
$wsKetamine = $wb.Worksheets.Item("Ketamine")
$wsKetamine.Range("D5").Value = 30
$wsKetamine.Range("D9").Value = 0.04

$excel.Calculate()
